$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (preserve rich-text runs by editing in place) ---
# A8 shared string: "Volume 32   Number  42" -> "...43"
$ws.Range("A8").Characters(21, 2).Text = "43"

# C9 shared string: "Report Covering the Week  10/13/2025  Through  10/19/2025"
#                 -> "Report Covering the Week  10/20/2025  Through  10/26/2025"
$ws.Range("C9").Characters(27, 10).Text = "10/20/2025"
$ws.Range("C9").Characters(48, 10).Text = "10/26/2025"

# --- Fix up cell styles/types for cells that change between numeric and
#     shared-text ("0" / "***.*") representations. Copying from a donor
#     cell that already has the desired style (and, where applicable,
#     shared-string value) brings across both in one step; numeric donors
#     are followed by an explicit .Value write below. ---
$ws.Range("C15").Copy($ws.Range("G15"))
$ws.Range("E15").Copy($ws.Range("H15"))
$ws.Range("D15").Copy($ws.Range("C20"))
$ws.Range("F15").Copy($ws.Range("D20"))
$ws.Range("N22").Copy($ws.Range("E20"))
$ws.Range("I15").Copy($ws.Range("C23"))
$ws.Range("C22").Copy($ws.Range("G27"))
$ws.Range("N23").Copy($ws.Range("H27"))
$ws.Range("C27").Copy($ws.Range("C28"))
$ws.Range("J15").Copy($ws.Range("D28"))
$ws.Range("K15").Copy($ws.Range("E28"))
$ws.Range("I20").Copy($ws.Range("G28"))
$ws.Range("L15").Copy($ws.Range("H28"))

# --- Numeric value updates for the crime-stats grid (rows 15-30) ---
$ws.Range("N15").Value = 200
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 6
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 6.666666666666
$ws.Range("I16").Value = 126
$ws.Range("J16").Value = 136
$ws.Range("K16").Value = -7.352941176470
$ws.Range("L16").Value = -12.5
$ws.Range("M16").Value = -35.714285714285
$ws.Range("N16").Value = 162.5
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 36.363636363636
$ws.Range("I17").Value = 225
$ws.Range("J17").Value = 202
$ws.Range("K17").Value = 11.386138613861
$ws.Range("L17").Value = 7.655502392344
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 675.862068965517
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -30.769230769230
$ws.Range("I18").Value = 91
$ws.Range("J18").Value = 98
$ws.Range("K18").Value = -7.142857142857
$ws.Range("L18").Value = -6.185567010309
$ws.Range("M18").Value = 2.247191011235
$ws.Range("N18").Value = 97.826086956521
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 303
$ws.Range("J19").Value = 297
$ws.Range("K19").Value = 2.020202020202
$ws.Range("L19").Value = 0.331125827814
$ws.Range("M19").Value = 37.10407239819
$ws.Range("N19").Value = 877.41935483871
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 28.571428571428
$ws.Range("L20").Value = -4.098360655737
$ws.Range("N20").Value = 143.75
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 70
$ws.Range("H21").Value = 2.857142857142
$ws.Range("I21").Value = 877
$ws.Range("J21").Value = 838
$ws.Range("K21").Value = 4.653937947494
$ws.Range("L21").Value = -1.238738738738
$ws.Range("M21").Value = 19.645293315143
$ws.Range("N21").Value = 323.671497584541
$ws.Range("D22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = -66.666666666666
$ws.Range("M22").Value = -64.285714285714
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -75
$ws.Range("I23").Value = 15
$ws.Range("J23").Value = 27
$ws.Range("K23").Value = -44.444444444444
$ws.Range("L23").Value = 15.384615384615
$ws.Range("M23").Value = 36.363636363636
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 23.529411764705
$ws.Range("F24").Value = 65
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = -19.753086419753
$ws.Range("I24").Value = 656
$ws.Range("J24").Value = 761
$ws.Range("K24").Value = -13.797634691195
$ws.Range("L24").Value = -16.326530612244
$ws.Range("M24").Value = 33.877551020408
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = -78.947368421052
$ws.Range("I25").Value = 132
$ws.Range("J25").Value = 247
$ws.Range("K25").Value = -46.558704453441
$ws.Range("L25").Value = -27.071823204419
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -62.5
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 308
$ws.Range("J26").Value = 284
$ws.Range("K26").Value = 8.450704225352
$ws.Range("L26").Value = -2.222222222222
$ws.Range("M26").Value = -24.509803921568
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 31
$ws.Range("K28").Value = 45.161290322580
$ws.Range("G29").Value = 2
$ws.Range("G30").Value = 1
